$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.091.79"
$ws.Range("E2").Value = "  -3.45%  "

$ws.Range("D3").Value = "3.521.38"
$ws.Range("E3").Value = "  -4.47%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'581.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").Value = "'175.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "

$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").Value = "3.515.42"
$ws.Range("E8").Value = "  -4.42%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -5.51%  "

$ws.Range("D11").Value = "'6.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.56%  "

$ws.Range("D12").Value = "'0.603"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("E13").Value = "  -4.92%  "

$ws.Range("D15").Value = "'675.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "4.087.03"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "3.521.77"
$ws.Range("E18").Value = "  -4.56%  "

$ws.Range("D19").Value = "69.089.92"
$ws.Range("E19").Value = "  -3.60%  "

$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").Value = "'17.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("D22").Value = "'11.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.17%  "

$ws.Range("D23").Value = "'0.911"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("D24").Value = "'16.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.47%  "

$ws.Range("D25").Value = "'98.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.43%  "

$ws.Range("E26").Value = "  -4.11%  "

$ws.Range("D27").Value = "'5.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("E28").Value = "  -5.76%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -6.69%  "

$ws.Range("E31").Value = "  -6.66%  "

$ws.Range("D32").Value = "'8.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.69%  "

$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.24%  "

$ws.Range("D34").Value = "'7.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "

$ws.Range("E35").Value = "  -5.15%  "

$ws.Range("D36").Value = "'579.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").Value = "'3.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -14.88%  "

$ws.Range("D38").Value = "'10.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.13%  "

$ws.Range("E39").Value = "  -3.32%  "

$ws.Range("E40").Value = "  -3.45%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("D43").Value = "'0.0441"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("E44").Value = "  -5.87%  "

$ws.Range("D45").Value = "3.437.23"
$ws.Range("E45").Value = "  -8.82%  "

$ws.Range("D46").Value = "'33.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.01%  "

$ws.Range("D47").Value = "0.0₃0711"
$ws.Range("E47").Value = "  -8.29%  "

$ws.Range("D48").Value = "'2.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("E49").Value = "  -6.37%  "

$ws.Range("D51").Value = "'131.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.94%  "
